# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 05:22"

# Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 886709
$ws.Range("C4").Value = 267
$ws.Range("E4").Value = 750544
$ws.Range("G4").Value = 7
$ws.Range("H4").Value = 50243

# Guatemala overtakes "Consejo Danes para los Refugiados" and "Sri Lanka"
# in the ranking (rows 112-114), shifting those two down a row.
$ws.Range("A112").Value = "Guatemala"
$ws.Range("B112").Value = 384
$ws.Range("C112").Value = 42
$ws.Range("D112").Value = 30
$ws.Range("E112").Value = 343
$ws.Range("F112").Value = 3
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 11

$ws.Range("A113").Value = "Consejo Danes para los Refugiados"
$ws.Range("B113").Value = 377
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 47
$ws.Range("E113").Value = 305
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 25

$ws.Range("A114").Value = "Sri Lanka"
$ws.Range("B114").Value = 368
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 107
$ws.Range("E114").Value = 254
$ws.Range("F114").Value = 2
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 7

# Birmania overtakes "Brunei" and "Gibraltar" in the ranking (rows 135-137),
# shifting those two down a row.
$ws.Range("A135").Value = "Birmania"
$ws.Range("B135").Value = 139
$ws.Range("C135").Value = 7
$ws.Range("D135").Value = 9
$ws.Range("E135").Value = 125
$ws.Range("F135").Value = 0
$ws.Range("G135").Value = 0
$ws.Range("H135").Value = 5

$ws.Range("A136").Value = "Brunei"
$ws.Range("B136").Value = 138
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 119
$ws.Range("E136").Value = 18
$ws.Range("F136").Value = 2
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 1

$ws.Range("A137").Value = "Gibraltar"
$ws.Range("B137").Value = 133
$ws.Range("C137").Value = 0
$ws.Range("D137").Value = 129
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 0
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 0

# Nepal (row 165) - minor update to active/recovered split
$ws.Range("D165").Value = 10
$ws.Range("E165").Value = 38
